$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Narrow column G (Recorded By) from width 31 to 13.
# Excel's ColumnWidth property is offset from the stored OOXML width by the
# default "extra padding" of 5/6 of a character, so subtract that to land on
# an exact stored width of 13.
$ws.Columns.Item(7).ColumnWidth = 12.166666666666666

# Replace every "Miss Dina Nasr, Administrator" entry in column G (Recorded By)
# with the academic year "2025/2026".
for ($r = 2; $r -le 673; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq "Miss Dina Nasr, Administrator") {
        $cell.Value = "2025/2026"
    }
}
